$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11 (pushes the old rows 11-12 down to
# 12-13, carrying their cell formatting/styles with them - e.g. the F-column
# hyperlink style survives the shift, including onto the new blank F11).
$ws.Rows("11:11").Insert()

# The whole scrape run stamps every data row with the new fetch timestamp.
$ts = "2025-11-07 18:24:50"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $ts
}

# Populate the newly inserted row with the newly scraped listing.
$ws.Cells.Item(11, 2).Value = "初回 ssss"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5429554"
$ws.Cells.Item(11, 7).Value = 18

# Only one brand-new hyperlink relationship gets appended (matching the
# source tool's insert-then-append behaviour where pre-existing hyperlink
# entries keep referencing their original rows rather than being
# renumbered); it lands on the row now holding the former last entry.
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://www.lancers.jp/work/detail/5428509")
# Re-assert the (already-shifted-down) Hyperlink style so the new
# relationship doesn't leave the cell on a freshly-minted duplicate style.
$ws.Cells.Item(13, 6).Style = "Hyperlink"
